# Weekly "Updated cryptos list" data refresh (GitHub Actions bot) applied to
# Sheet1. Only the Price (D) and Volume(1h) (E) columns move for most rows;
# rows 38/39 (WhiteBITCoin <-> PolygonEcosystemToken) additionally swapped
# rank position, so their Coin name (B) and Link (C) are rewritten too.
#
# Column D is stored as literal text (e.g. "76.412.52", "3.045.63" are
# thousands-grouped price strings, not numbers). A handful of the refreshed
# prices (e.g. "624.54", "1.00", "3.90") would otherwise parse as genuine
# numbers, and a plain .Value assignment lets Excel silently coerce them to
# numeric cells (dropping significant trailing zeros, e.g. "3.90" -> 3.9).
# Writing those through .Formula with a leading "'" forces a literal-text
# entry instead, so the stored text matches the source data exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '76.389.18'
$ws.Range("E2").Value = '  +0.35%  '

$ws.Range("D3").Value = '3.049.00'
$ws.Range("E3").Value = '  +3.67%  '

$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("E5").Value = '  -1.37%  '

$ws.Range("D6").Formula = "'" + '624.54'
$ws.Range("E6").Value = '  +4.45%  '

$ws.Range("D7").Formula = "'" + '1.00'
$ws.Range("E7").Value = '  +0.09%  '

$ws.Range("D8").Formula = "'" + '0.549'
$ws.Range("E8").Value = '  -0.56%  '

$ws.Range("D9").Formula = "'" + '0.206'
$ws.Range("E9").Value = '  +4.07%  '

$ws.Range("D10").Value = '3.044.90'
$ws.Range("E10").Value = '  +3.88%  '

$ws.Range("D11").Formula = "'" + '0.442'
$ws.Range("E11").Value = '  +1.09%  '

$ws.Range("E13").Value = '  +5.68%  '

$ws.Range("D14").Value = '3.605.88'
$ws.Range("E14").Value = '  +3.92%  '

$ws.Range("D15").Formula = "'" + '29.11'
$ws.Range("E15").Value = '  +3.64%  '

$ws.Range("D16").Value = '76.301.45'
$ws.Range("E16").Value = '  +0.45%  '

$ws.Range("E17").Value = '  +1.50%  '

$ws.Range("D18").Value = '3.046.61'
$ws.Range("E18").Value = '  +3.67%  '

$ws.Range("D19").Formula = "'" + '13.55'
$ws.Range("E19").Value = '  +2.51%  '

$ws.Range("E20").Value = '  +1.41%  '

$ws.Range("D21").Formula = "'" + '375.27'
$ws.Range("E21").Value = '  +0.54%  '

$ws.Range("D22").Formula = "'" + '2.31'
$ws.Range("E22").Value = '  +0.30%  '

$ws.Range("D23").Formula = "'" + '4.37'
$ws.Range("E23").Value = '  +1.66%  '

$ws.Range("D24").Formula = "'" + '73.32'
$ws.Range("E24").Value = '  +2.21%  '

$ws.Range("E26").Value = '  -0.17%  '

$ws.Range("E27").Value = '  +1.20%  '

$ws.Range("D28").Formula = "'" + '9.80'
$ws.Range("E28").Value = '  +0.91%  '

$ws.Range("E29").Value = '  +0.73%  '

$ws.Range("D30").Formula = "'" + '1.00'
$ws.Range("E30").Value = '  +0.04%  '

$ws.Range("D31").Formula = "'" + '8.30'
$ws.Range("E31").Value = '  +6.25%  '

$ws.Range("E32").Value = '  +0.92%  '

$ws.Range("E33").Value = '  +5.69%  '

$ws.Range("D34").Formula = "'" + '492.71'
$ws.Range("E34").Value = '  -2.36%  '

$ws.Range("E35").Value = '  -0.04%  '

$ws.Range("D36").Formula = "'" + '20.68'
$ws.Range("E36").Value = '  +1.82%  '

$ws.Range("D37").Formula = "'" + '162.81'
$ws.Range("E37").Value = '  -0.62%  '

$ws.Range("B38").Value = 'PolygonEcosystemToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D38").Formula = "'" + '0.384'
$ws.Range("E38").Value = '  +3.52%  '

$ws.Range("B39").Value = 'WhiteBITCoin'
$ws.Range("C39").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D39").Formula = "'" + '20.04'
$ws.Range("E39").Value = '  +2.08%  '

$ws.Range("E40").Value = '  +2.68%  '

$ws.Range("D41").Formula = "'" + '190.81'
$ws.Range("E41").Value = '  +5.03%  '

$ws.Range("E42").Value = '  -6.28%  '

$ws.Range("D44").Formula = "'" + '0.812'
$ws.Range("E44").Value = '  +22.85%  '

$ws.Range("E45").Value = '  +2.77%  '

$ws.Range("E46").Value = '  +5.43%  '

$ws.Range("D47").Formula = "'" + '42.06'
$ws.Range("E47").Value = '  +4.76%  '

$ws.Range("E48").Value = '  -0.80%  '

$ws.Range("D49").Formula = "'" + '2.46'
$ws.Range("E49").Value = '  +4.73%  '

$ws.Range("E50").Value = '  +4.28%  '

$ws.Range("D51").Formula = "'" + '3.90'
$ws.Range("E51").Value = '  +3.87%  '
